$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.999.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.224.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.88%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -6.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -10.55%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.225.79"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.459"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.63"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.776.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -8.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.215.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000159"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.061.84"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "362.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.520"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.355.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.171"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0977"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -11.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.08"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.93"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.24"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.72"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.39"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.76%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.57"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -15.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0704"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.53%  "
$ws.Range("B42").Value = "RenzoRestakedETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.256.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.723"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.02"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.51"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.290.78"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.26"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.92%  "
